$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(3)

# Grow the textbox height to fit the extra bullet content (width/position unchanged).
# Shape.Height is expressed in points; 10813312 EMU wide box keeps its cx/off, only cy grows
# from 2554545 EMU to 2862322 EMU == 225.37973 pt (912775/4048 scaled by 1/12700).
$sh.Height = 225.37973

$tr = $sh.TextFrame.TextRange

# Work from the last paragraph upward so earlier paragraph indices stay valid
# while paragraphs are inserted/removed below them.

# Paragraph 8: "No fee paid to the central Carpentries organization (approx. $XXX HOW MUCH WAS IT?)"
# -> drop the parenthetical and split the sentence into two runs.
$para8 = $tr.Paragraphs(8, 1)
$para8.Text = "No fee paid to the central Carpentries organization"
$firstPart = "No fee paid to the central "
$secondPart = "Carpentries organization"
$secondRun = $para8.Characters($firstPart.Length + 1, $secondPart.Length)
$secondRun.Text = $secondPart

# Paragraph 7: empty placeholder bullet right above "No fee..." -> remove it.
[void]$tr.Paragraphs(7, 1).Delete()

# Paragraph 6 ("instructor buy-in") gains a new sibling sub-bullet after it.
[void]$tr.Paragraphs(6, 1).InsertAfter([char]13 + "curriculum ")

# Paragraph 2: empty placeholder sub-bullet right under the opening bullet -> remove it.
[void]$tr.Paragraphs(2, 1).Delete()

# Paragraph 1: reword the opening bullet and add a brand-new bullet right after it.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "NNLM sponsored Carpentries Instructor training building its network of certified instructors and capacity to self-organize carpentries workshops. "
[void]$para1.InsertAfter([char]13 + "HSHSL built its capacity with 3 in-house staff gaining LC instructor certification.")
